$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33, pushing the existing row 33 (and below) down to 34.
$ws.Rows("33:33").Insert()

# Populate the new row 33 with the latest weekly price observation.
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44706
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100112010
$ws.Range("G33").Value = "Achicoria"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 11000
$ws.Range("M33").Value = 11000
$ws.Range("N33").Value = "$/caja 18 unidades"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 611
$ws.Range("Q33").Value = 18
$ws.Range("R33").Value = "Hortaliza"
